$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data-driven update of rows 2-51 per the diff (price/volume refresh, timestamp bump, and a few name/link row shifts).
$changes = @(
    @{Row=2; D='242.19'; E='-1.45%'; G='2'};
    @{Row=3; D='29.16'; E='12.04%'; G='2'};
    @{Row=4; D='5.093'; E='-1.29%'; G='2'};
    @{Row=5; D='0.05641'; E='0.95%'; G='2'};
    @{Row=6; D='6.504'; E='0.28%'; G='2'};
    @{Row=7; D='0.8201'; E='0.79%'; G='2'};
    @{Row=8; D='0.8542'; E='1.64%'; G='2'};
    @{Row=9; D='0.1331'; E='0.48%'; G='2'};
    @{Row=10; B='BitrueCoin'; C='https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'; D='0.02864'; E='0.43%'; G='2'};
    @{Row=11; B='BitMartToken'; C='https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'; D='0.09371'; E='-0.06%'; G='2'};
    @{Row=12; B='BitForexToken'; C='https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'; D='0.001516'; E='0.35%'; G='2'};
    @{Row=13; B='CoinExToken'; C='https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'; D='0.04147'; E='-10.27%'; G='2'};
    @{Row=14; B='One'; C='https://coinranking.com/coin/6Lga5NiXX3rT+one-one'; D='0.0006023'; E='-93.95%'; G='2'};
    @{Row=15; B='TigerCash'; C='https://coinranking.com/coin/6hIn06L2+tigercash-tch'; D='0.006218'; E='0.83%'; G='2'};
    @{Row=16; B='LEO'; C='https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'; D='3.526'; E='-2.59%'; G='2'};
    @{Row=17; B='GateToken'; C='https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'; D='3.013'; E='-0.71%'; G='2'};
    @{Row=18; B='BTSEToken'; C='https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'; D='2.218'; E='1.59%'; G='2'};
    @{Row=19; B='BitpandaEcosystemToken'; C='https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'; D='0.3152'; E='1.26%'; G='2'};
    @{Row=20; B='MandalaExchangeToken'; C='https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'; D='0.06873'; E='-1.59%'; G='2'};
    @{Row=21; D='0.03219'; E='3.53%'; G='2'};
    @{Row=22; D='0.1273'; E='-2.05%'; G='2'};
    @{Row=23; D='3.619'; E='-3.69%'; G='2'};
    @{Row=24; E='0.10%'; G='2'};
    @{Row=25; D='0.001213'; E='-2.79%'; G='2'};
    @{Row=26; D='0.004455'; E='-1.63%'; G='2'};
    @{Row=27; D='0.0001181'; E='23.06%'; G='2'};
    @{Row=28; E='-15.26%'; G='2'};
    @{Row=29; G='2'};
    @{Row=30; G='2'};
    @{Row=31; G='2'};
    @{Row=32; G='2'};
    @{Row=33; G='2'};
    @{Row=34; G='2'};
    @{Row=35; G='2'};
    @{Row=36; G='2'};
    @{Row=37; G='2'};
    @{Row=38; G='2'};
    @{Row=39; G='2'};
    @{Row=40; D='0.03707'; E='1.83%'; G='2'};
    @{Row=41; D='0.005793'; E='-6.35%'; G='2'};
    @{Row=42; D='0.1053'; E='0.30%'; G='2'};
    @{Row=43; E='-8.71%'; G='2'};
    @{Row=44; D='0.009449'; E='4.22%'; G='2'};
    @{Row=45; D='0.00005107'; E='-4.63%'; G='2'};
    @{Row=46; D='0.00000000750'; E='0.11%'; G='2'};
    @{Row=47; D='0.1201'; E='10.20%'; G='2'};
    @{Row=48; D='0.002525'; E='-1.08%'; G='2'};
    @{Row=49; E='0.11%'; G='2'};
    @{Row=50; E='0.11%'; G='2'};
    @{Row=51; G='2'};
)

foreach ($item in $changes) {
    $row = $item.Row

    if ($item.ContainsKey('B')) {
        $ws.Cells.Item($row, 2).Value = $item.B
    }
    if ($item.ContainsKey('C')) {
        $ws.Cells.Item($row, 3).Value = $item.C
    }
    if ($item.ContainsKey('D')) {
        $cell = $ws.Cells.Item($row, 4)
        $cell.Value = "'" + $item.D
        $cell.Style = "Normal"
    }
    if ($item.ContainsKey('E')) {
        $cell = $ws.Cells.Item($row, 5)
        $cell.Value = "'" + $item.E
        $cell.Style = "Normal"
    }
    if ($item.ContainsKey('G')) {
        $cell = $ws.Cells.Item($row, 7)
        $cell.Value = "'" + $item.G
        $cell.Style = "Normal"
    }
}
